# Applies updated profit/price figures to each job sheet as produced by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1970.7778
$ws.Range("I12").Value = 2574.5
$ws.Range("K12").Value = 2574.5
$ws.Range("M12").Value = -2404.5
$ws.Range("H19").Value = 11906292
$ws.Range("I19").Value = 1146.0714
$ws.Range("K19").Value = 1146.0714
$ws.Range("M19").Value = -971.0714
$ws.Range("H29").Value = 4019.8
$ws.Range("I29").Value = 299.5
$ws.Range("K29").Value = 898.5
$ws.Range("M29").Value = -617.5
$ws.Range("H64").Value = 4857.4
$ws.Range("I64").Value = 4802
$ws.Range("J64").Value = 4871.25
$ws.Range("K64").Value = 4802
$ws.Range("L64").Value = 4871.25
$ws.Range("M64").Value = -4554
$ws.Range("N64").Value = -5367.25
$ws.Range("H67").Value = 4857.4
$ws.Range("I67").Value = 4802
$ws.Range("J67").Value = 4871.25
$ws.Range("K67").Value = 4802
$ws.Range("L67").Value = 4871.25
$ws.Range("M67").Value = -3944
$ws.Range("N67").Value = -6587.25
$ws.Range("H76").Value = 55616616
$ws.Range("I76").Value = 170563.83
$ws.Range("J76").Value = 83339640
$ws.Range("K76").Value = 170563.83
$ws.Range("L76").Value = 83339640
$ws.Range("M76").Value = -170248.83
$ws.Range("N76").Value = -83340270
$ws.Range("H79").Value = 55616616
$ws.Range("I79").Value = 170563.83
$ws.Range("J79").Value = 83339640
$ws.Range("K79").Value = 170563.83
$ws.Range("L79").Value = 83339640
$ws.Range("M79").Value = -169471.83
$ws.Range("N79").Value = -83341824
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
$ws.Range("H87").Value = 74535
$ws.Range("J87").Value = 74812.19
$ws.Range("L87").Value = 74812.19
$ws.Range("N87").Value = -77308.19
$ws.Range("H90").Value = 74535
$ws.Range("J90").Value = 74812.19
$ws.Range("L90").Value = 224436.57
$ws.Range("N90").Value = -236916.57
$ws.Range("H96").Value = 681.3333
$ws.Range("I96").Value = 419
$ws.Range("J96").Value = 1993
$ws.Range("K96").Value = 1257
$ws.Range("L96").Value = 5979
$ws.Range("M96").Value = 116
$ws.Range("N96").Value = -8725
$ws.Range("H97").Value = 1949
$ws.Range("I97").Value = 499
$ws.Range("K97").Value = 1497
$ws.Range("M97").Value = -1001
$ws.Range("H107").Value = 26891.18
$ws.Range("I107").Value = 37057
$ws.Range("J107").Value = 1014.5455
$ws.Range("K107").Value = 37057
$ws.Range("L107").Value = 1014.5455
$ws.Range("M107").Value = -35137
$ws.Range("N107").Value = -4854.5455
$ws.Range("H111").Value = 73686.92999999999
$ws.Range("I111").Value = 201645
$ws.Range("J111").Value = 2599.111
$ws.Range("K111").Value = 604935
$ws.Range("L111").Value = 7797.333
$ws.Range("M111").Value = -601868
$ws.Range("N111").Value = -13931.333
$ws.Range("H113").Value = 9900
$ws.Range("I113").Value = 9866.666999999999
$ws.Range("K113").Value = 9866.666999999999
$ws.Range("M113").Value = -6612.666999999999
$ws.Range("H118").Value = 38321.8
$ws.Range("I118").Value = 38321.8
$ws.Range("K118").Value = 114965.4
$ws.Range("M118").Value = -113308.4
$ws.Range("H127").Value = 1311.1818
$ws.Range("J127").Value = 2749.75
$ws.Range("L127").Value = 8249.25
$ws.Range("N127").Value = -18169.25
$ws.Range("H129").Value = 2342.3416
$ws.Range("I129").Value = 1024.4286
$ws.Range("J129").Value = 2613.6765
$ws.Range("K129").Value = 3073.2858
$ws.Range("L129").Value = 7841.029500000001
$ws.Range("M129").Value = 1926.7142
$ws.Range("N129").Value = -17841.0295
$ws.Range("H132").Value = 5018.171
$ws.Range("I132").Value = 5124.3076
$ws.Range("J132").Value = 2948.5
$ws.Range("K132").Value = 15372.9228
$ws.Range("L132").Value = 8845.5
$ws.Range("M132").Value = -12842.9228
$ws.Range("N132").Value = -13905.5
$ws.Range("H135").Value = 3820.96
$ws.Range("I135").Value = 3820.96
$ws.Range("K135").Value = 34388.64
$ws.Range("M135").Value = -31853.64
$ws.Range("H137").Value = 22842.857
$ws.Range("I137").Value = 8088.25
$ws.Range("J137").Value = 26314.53
$ws.Range("K137").Value = 24264.75
$ws.Range("L137").Value = 78943.59
$ws.Range("M137").Value = -21714.75
$ws.Range("N137").Value = -84043.59
$ws.Range("H138").Value = 6494.44
$ws.Range("I138").Value = 2086.35
$ws.Range("J138").Value = 8097.382
$ws.Range("K138").Value = 6259.049999999999
$ws.Range("L138").Value = 24292.146
$ws.Range("M138").Value = -1119.049999999999
$ws.Range("N138").Value = -34572.146

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3936.353
$ws.Range("I32").Value = 3892.0652
$ws.Range("K32").Value = 3892.0652
$ws.Range("M32").Value = -3605.0652
$ws.Range("H45").Value = 80166.92
$ws.Range("I45").Value = 114473.22
$ws.Range("J45").Value = 2977.75
$ws.Range("K45").Value = 114473.22
$ws.Range("L45").Value = 2977.75
$ws.Range("M45").Value = -114096.22
$ws.Range("N45").Value = -3731.75
$ws.Range("H61").Value = 4690.7
$ws.Range("I61").Value = 4363.375
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 4363.375
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = -4151.375
$ws.Range("N61").Value = -6424
$ws.Range("H63").Value = 7589.9
$ws.Range("J63").Value = 9999.857
$ws.Range("L63").Value = 9999.857
$ws.Range("N63").Value = -11371.857
$ws.Range("H66").Value = 7589.9
$ws.Range("J66").Value = 9999.857
$ws.Range("L66").Value = 49999.285
$ws.Range("N66").Value = -56863.285
$ws.Range("H74").Value = 512291.34
$ws.Range("I74").Value = 558099.9399999999
$ws.Range("K74").Value = 558099.9399999999
$ws.Range("M74").Value = -557225.9399999999
$ws.Range("H77").Value = 512291.34
$ws.Range("I77").Value = 558099.9399999999
$ws.Range("K77").Value = 2790499.7
$ws.Range("M77").Value = -2786131.7
$ws.Range("H122").Value = 58827176
$ws.Range("I122").Value = 125001480
$ws.Range("K122").Value = 375004440
$ws.Range("M122").Value = -375001990
$ws.Range("H124").Value = 64123
$ws.Range("J124").Value = 64123
$ws.Range("L124").Value = 64123
$ws.Range("N124").Value = -73943
$ws.Range("H132").Value = 180697.62
$ws.Range("I132").Value = 235729.2
$ws.Range("K132").Value = 707187.6000000001
$ws.Range("M132").Value = -704657.6000000001
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140
$ws.Range("H136").Value = 4690.7
$ws.Range("I136").Value = 4363.375
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 13090.125
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -10540.125
$ws.Range("N136").Value = -23100

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
$ws.Range("H20").Value = 2058.2964
$ws.Range("I20").Value = 1984.625
$ws.Range("K20").Value = 1984.625
$ws.Range("M20").Value = -1737.625
$ws.Range("H29").Value = 18703.2
$ws.Range("I29").Value = 4338.6665
$ws.Range("K29").Value = 4338.6665
$ws.Range("M29").Value = -4049.6665
$ws.Range("H36").Value = 22599.8
$ws.Range("J36").Value = 1000
$ws.Range("L36").Value = 1000
$ws.Range("N36").Value = -2068
$ws.Range("H82").Value = 55427.273
$ws.Range("J82").Value = 95998
$ws.Range("L82").Value = 95998
$ws.Range("N82").Value = -96764
$ws.Range("H85").Value = 55427.273
$ws.Range("J85").Value = 95998
$ws.Range("L85").Value = 95998
$ws.Range("N85").Value = -98650
$ws.Range("H99").Value = 1030.0588
$ws.Range("I99").Value = 900.7857
$ws.Range("K99").Value = 900.7857
$ws.Range("M99").Value = 597.2143
$ws.Range("H134").Value = 30871.422
$ws.Range("I134").Value = 1669.697
$ws.Range("K134").Value = 5009.090999999999
$ws.Range("M134").Value = -2474.090999999999

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2277.8635
$ws.Range("I16").Value = 1969.125
$ws.Range("J16").Value = 3101.1667
$ws.Range("K16").Value = 1969.125
$ws.Range("L16").Value = 3101.1667
$ws.Range("M16").Value = -1682.125
$ws.Range("N16").Value = -3675.1667
$ws.Range("H22").Value = 289.8
$ws.Range("I22").Value = 237.25
$ws.Range("K22").Value = 237.25
$ws.Range("M22").Value = 112.75
$ws.Range("H31").Value = 406143.88
$ws.Range("I31").Value = 943028.25
$ws.Range("J31").Value = 19587.16
$ws.Range("K31").Value = 943028.25
$ws.Range("L31").Value = 19587.16
$ws.Range("M31").Value = -942733.25
$ws.Range("N31").Value = -20177.16
$ws.Range("H34").Value = 406143.88
$ws.Range("I34").Value = 943028.25
$ws.Range("J34").Value = 19587.16
$ws.Range("K34").Value = 943028.25
$ws.Range("L34").Value = 19587.16
$ws.Range("M34").Value = -942826.25
$ws.Range("N34").Value = -19991.16
$ws.Range("H58").Value = 8171.75
$ws.Range("I58").Value = 2719.7334
$ws.Range("J58").Value = 14462.538
$ws.Range("K58").Value = 2719.7334
$ws.Range("L58").Value = 14462.538
$ws.Range("M58").Value = -2516.7334
$ws.Range("N58").Value = -14868.538
$ws.Range("H93").Value = 39754.5
$ws.Range("I93").Value = 39754.5
$ws.Range("K93").Value = 39754.5
$ws.Range("M93").Value = -37882.5
$ws.Range("H99").Value = 5414.0586
$ws.Range("I99").Value = 4724
$ws.Range("J99").Value = 6399.857
$ws.Range("K99").Value = 4724
$ws.Range("L99").Value = 6399.857
$ws.Range("M99").Value = -3226
$ws.Range("N99").Value = -9395.857
$ws.Range("H105").Value = 1678.1
$ws.Range("I105").Value = 1314.8334
$ws.Range("K105").Value = 1314.8334
$ws.Range("M105").Value = 432.1666
$ws.Range("H113").Value = 2277.8635
$ws.Range("I113").Value = 1969.125
$ws.Range("J113").Value = 3101.1667
$ws.Range("K113").Value = 1969.125
$ws.Range("L113").Value = 3101.1667
$ws.Range("M113").Value = 200.875
$ws.Range("N113").Value = -7441.1667
$ws.Range("H122").Value = 40034.223
$ws.Range("I122").Value = 85295.75
$ws.Range("K122").Value = 255887.25
$ws.Range("M122").Value = -253437.25
$ws.Range("H126").Value = 5414.0586
$ws.Range("I126").Value = 4724
$ws.Range("J126").Value = 6399.857
$ws.Range("K126").Value = 14172
$ws.Range("L126").Value = 19199.571
$ws.Range("M126").Value = -11702
$ws.Range("N126").Value = -24139.571
$ws.Range("H132").Value = 2915.1538
$ws.Range("I132").Value = 1945.2273
$ws.Range("K132").Value = 5835.6819
$ws.Range("M132").Value = -3305.6819
$ws.Range("H134").Value = 243751.28
$ws.Range("I134").Value = 2618.1
$ws.Range("J134").Value = 846584.25
$ws.Range("K134").Value = 7854.299999999999
$ws.Range("L134").Value = 2539752.75
$ws.Range("M134").Value = -5319.299999999999
$ws.Range("N134").Value = -2544822.75
$ws.Range("H136").Value = 8171.75
$ws.Range("I136").Value = 2719.7334
$ws.Range("J136").Value = 14462.538
$ws.Range("K136").Value = 8159.2002
$ws.Range("L136").Value = 43387.614
$ws.Range("M136").Value = -5609.2002
$ws.Range("N136").Value = -48487.614

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 4050.6
$ws.Range("J2").Value = 5007.5
$ws.Range("L2").Value = 30045
$ws.Range("N2").Value = -30271
$ws.Range("H5").Value = 1214691.8
$ws.Range("J5").Value = 13561.75
$ws.Range("L5").Value = 40685.25
$ws.Range("N5").Value = -40909.25
$ws.Range("H14").Value = 6935.5
$ws.Range("I14").Value = 6935.5
$ws.Range("K14").Value = 20806.5
$ws.Range("M14").Value = -20633.5
$ws.Range("H50").Value = 543.9375
$ws.Range("I50").Value = 341.33334
$ws.Range("J50").Value = 590.6923
$ws.Range("K50").Value = 1024.00002
$ws.Range("L50").Value = 1772.0769
$ws.Range("M50").Value = -543.0000199999999
$ws.Range("N50").Value = -2734.0769
$ws.Range("H53").Value = 543.9375
$ws.Range("I53").Value = 341.33334
$ws.Range("J53").Value = 590.6923
$ws.Range("K53").Value = 1024.00002
$ws.Range("L53").Value = 1772.0769
$ws.Range("M53").Value = -543.0000199999999
$ws.Range("N53").Value = -2734.0769
$ws.Range("H68").Value = 1926.1389
$ws.Range("I68").Value = 1430.5
$ws.Range("J68").Value = 2067.75
$ws.Range("K68").Value = 4291.5
$ws.Range("L68").Value = 6203.25
$ws.Range("M68").Value = -3480.5
$ws.Range("N68").Value = -7825.25
$ws.Range("H71").Value = 1926.1389
$ws.Range("I71").Value = 1430.5
$ws.Range("J71").Value = 2067.75
$ws.Range("K71").Value = 12874.5
$ws.Range("L71").Value = 18609.75
$ws.Range("M71").Value = -8818.5
$ws.Range("N71").Value = -26721.75
$ws.Range("H92").Value = 588.1
$ws.Range("J92").Value = 603.44446
$ws.Range("L92").Value = 1810.33338
$ws.Range("N92").Value = -4306.33338
$ws.Range("H101").Value = 15000
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 45000
$ws.Range("N101").Value = -49868
$ws.Range("H132").Value = 3052090.5
$ws.Range("J132").Value = 29247
$ws.Range("L132").Value = 263223
$ws.Range("N132").Value = -268283
$ws.Range("H135").Value = 1214691.8
$ws.Range("J135").Value = 13561.75
$ws.Range("L135").Value = 122055.75
$ws.Range("N135").Value = -127125.75
$ws.Range("H137").Value = 1996.9
$ws.Range("I137").Value = 1996.9
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5990.700000000001
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -890.7000000000007
$ws.Range("N137").Value = $null

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4272.6875
$ws.Range("I70").Value = 4096.222
$ws.Range("J70").Value = 4499.5713
$ws.Range("K70").Value = 4096.222
$ws.Range("L70").Value = 4499.5713
$ws.Range("M70").Value = -3826.222
$ws.Range("N70").Value = -5039.5713
$ws.Range("H73").Value = 4272.6875
$ws.Range("I73").Value = 4096.222
$ws.Range("J73").Value = 4499.5713
$ws.Range("K73").Value = 4096.222
$ws.Range("L73").Value = 4499.5713
$ws.Range("M73").Value = -3160.222
$ws.Range("N73").Value = -6371.5713
$ws.Range("H122").Value = 465633.78
$ws.Range("I122").Value = 556655.0600000001
$ws.Range("J122").Value = 10527.5
$ws.Range("K122").Value = 1669965.18
$ws.Range("L122").Value = 31582.5
$ws.Range("M122").Value = -1667515.18
$ws.Range("N122").Value = -36482.5
$ws.Range("H132").Value = 17576.422
$ws.Range("I132").Value = 1921.8545
$ws.Range("J132").Value = 113243.22
$ws.Range("K132").Value = 5765.5635
$ws.Range("L132").Value = 339729.66
$ws.Range("M132").Value = -3235.5635
$ws.Range("N132").Value = -344789.66
$ws.Range("H138").Value = 48730
$ws.Range("J138").Value = 48730
$ws.Range("L138").Value = 48730
$ws.Range("N138").Value = -59010

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 536243.4399999999
$ws.Range("I7").Value = 718972.9
$ws.Range("J7").Value = 24601
$ws.Range("K7").Value = 718972.9
$ws.Range("L7").Value = 24601
$ws.Range("M7").Value = -718860.9
$ws.Range("N7").Value = -24825
$ws.Range("H16").Value = 1910.2084
$ws.Range("I16").Value = 1907.75
$ws.Range("K16").Value = 1907.75
$ws.Range("M16").Value = -1737.75
$ws.Range("H20").Value = 50000
$ws.Range("J20").Value = 50000
$ws.Range("L20").Value = 50000
$ws.Range("N20").Value = -50452
$ws.Range("H61").Value = 4163
$ws.Range("I61").Value = 2801.2778
$ws.Range("K61").Value = 2801.2778
$ws.Range("M61").Value = -2599.2778
$ws.Range("H93").Value = 100002260
$ws.Range("I93").Value = 111113550
$ws.Range("K93").Value = 111113550
$ws.Range("M93").Value = -111112302
$ws.Range("H113").Value = 4163
$ws.Range("I113").Value = 2801.2778
$ws.Range("K113").Value = 2801.2778
$ws.Range("M113").Value = -631.2777999999998
$ws.Range("H122").Value = 793700.7
$ws.Range("I122").Value = 5300.3335
$ws.Range("J122").Value = 1385001
$ws.Range("K122").Value = 15901.0005
$ws.Range("L122").Value = 4155003
$ws.Range("M122").Value = -13451.0005
$ws.Range("N122").Value = -4159903
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null
$ws.Range("H126").Value = 536243.4399999999
$ws.Range("I126").Value = 718972.9
$ws.Range("J126").Value = 24601
$ws.Range("K126").Value = 2156918.7
$ws.Range("L126").Value = 73803
$ws.Range("M126").Value = -2154448.7
$ws.Range("N126").Value = -78743
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = $null
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null
$ws.Range("H132").Value = 6126.3687
$ws.Range("I132").Value = 5208.4165
$ws.Range("J132").Value = 7700
$ws.Range("K132").Value = 15625.2495
$ws.Range("L132").Value = 23100
$ws.Range("M132").Value = -13095.2495
$ws.Range("N132").Value = -28160

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7272.909
$ws.Range("I62").Value = 7250
$ws.Range("J62").Value = 7286
$ws.Range("K62").Value = 7250
$ws.Range("L62").Value = 7286
$ws.Range("M62").Value = -6626
$ws.Range("N62").Value = -8534
$ws.Range("H65").Value = 7272.909
$ws.Range("I65").Value = 7250
$ws.Range("J65").Value = 7286
$ws.Range("K65").Value = 36250
$ws.Range("L65").Value = 36430
$ws.Range("M65").Value = -33130
$ws.Range("N65").Value = -42670
$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524
$ws.Range("H109").Value = 68747
$ws.Range("I109").Value = 50000
$ws.Range("J109").Value = 74996
$ws.Range("K109").Value = 50000
$ws.Range("L109").Value = 74996
$ws.Range("M109").Value = -48613
$ws.Range("N109").Value = -77770
$ws.Range("H113").Value = 832.4583
$ws.Range("I113").Value = 986
$ws.Range("J113").Value = 371.83334
$ws.Range("K113").Value = 2958
$ws.Range("L113").Value = 1115.50002
$ws.Range("M113").Value = -788
$ws.Range("N113").Value = -5455.500019999999
$ws.Range("H122").Value = 4388.76
$ws.Range("I122").Value = 3107.5
$ws.Range("J122").Value = 6666.5557
$ws.Range("K122").Value = 9322.5
$ws.Range("L122").Value = 19999.6671
$ws.Range("M122").Value = -6872.5
$ws.Range("N122").Value = -24899.6671
$ws.Range("H132").Value = 40651.035
$ws.Range("I132").Value = 4120.524
$ws.Range("K132").Value = 12361.572
$ws.Range("M132").Value = -9831.572
$ws.Range("H136").Value = 438296.47
$ws.Range("I136").Value = 482252.25
$ws.Range("K136").Value = 1446756.75
$ws.Range("M136").Value = -1444206.75
